$wb = $excel.ActiveWorkbook

# --- Sheet "2025" (sheet1.xml) ---
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("A2").Value = 0
$ws2025.Range("B2").Value = 4876.512737957041
$ws2025.Range("E2").Value = 212279.6285102097
$ws2025.Range("G2").Value = 56671.47998863283
$ws2025.Range("I2").Value = 109653.27140292
$ws2025.Range("L2").Value = 428573.6569267786
$ws2025.Range("M2").Value = 80341.74698896499
$ws2025.Range("N2").Value = 44942.73878197405
$ws2025.Range("O2").Value = 50827.92708717833

# --- Sheet "2030" (sheet2.xml) ---
$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Value = 8893.010627692065
$ws2030.Range("B2").Value = 34788.11346181522
$ws2030.Range("E2").Value = 214922.2077203031
$ws2030.Range("I2").Value = 229183.4540128457
$ws2030.Range("L2").Value = 100980.5016519909
$ws2030.Range("M2").Value = 92065.42697618291
$ws2030.Range("N2").Value = 49188.66652121952
$ws2030.Range("O2").Value = 39944.84862230343
